# Horarios 141 update - 2026-01-18 08:38:24 scrape
# Inserts newly scraped rows into each sheet's table (sorted by arrival time),
# refreshes the "Última actualización" / "Total filas" header lines.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 08:38:24"
$ws1.Range("A3").Value = "Total filas: 75"

# New row inserted before the current row 41
$ws1.Rows.Item(41).Insert()
$ws1.Range("A41").Value = "08:38:24"
$ws1.Range("B41").Value = "08:40"
$ws1.Range("C41").Value = "10_OLMOS"
$ws1.Range("D41").Value = 2
$ws1.Range("E41").Value = "LP1912"

# New row inserted before the current row 61
$ws1.Rows.Item(61).Insert()
$ws1.Range("A61").Value = "08:38:24"
$ws1.Range("B61").Value = "09:17"
$ws1.Range("C61").Value = "27_EL RETIRO"
$ws1.Range("D61").Value = 39
$ws1.Range("E61").Value = "LP1912"

# New row inserted before the current row 69
$ws1.Rows.Item(69).Insert()
$ws1.Range("A69").Value = "08:38:24"
$ws1.Range("B69").Value = "09:34"
$ws1.Range("C69").Value = "23_HERNANDEZ"
$ws1.Range("D69").Value = 56
$ws1.Range("E69").Value = "LP1912"

# New row appended at the end (row 80)
$ws1.Range("A80").Value = "08:38:24"
$ws1.Range("B80").Value = "10:29"
$ws1.Range("C80").Value = "15_ABASTO"
$ws1.Range("D80").Value = 111
$ws1.Range("E80").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 08:38:24"

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 08:38:24"
$ws3.Range("A3").Value = "Total filas: 12"

# New row inserted before the current row 10
$ws3.Rows.Item(10).Insert()
$ws3.Range("A10").Value = "08:38:24"
$ws3.Range("B10").Value = "08:38"
$ws3.Range("C10").Value = "215A_LA PLATA"
$ws3.Range("D10").Value = 0
$ws3.Range("E10").Value = "L6173"

# New rows appended at the end (rows 16 and 17)
$ws3.Range("A16").Value = "08:38:24"
$ws3.Range("B16").Value = "10:13"
$ws3.Range("C16").Value = "215C_LA PLATA"
$ws3.Range("D16").Value = 95
$ws3.Range("E16").Value = "L6203"

$ws3.Range("A17").Value = "08:38:24"
$ws3.Range("B17").Value = "10:30"
$ws3.Range("C17").Value = "215B_LP-P MOR-1 Y 57"
$ws3.Range("D17").Value = 112
$ws3.Range("E17").Value = "L6173"
